# feat: add 2022-Q4 data
#
# Inserts a new "2022-Q4" sheet (positioned right after "总计" and before
# "2022-Q3"), populates it with the new quarter's fund-holdings data, and
# updates the "总计" (totals) sheet so the new quarter's summary row is
# inserted at the top of the data (pushing 2022-Q3 / 2022-Q1 down).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Update the "总计" summary sheet: insert a new 2022-Q4 row above the
#    existing 2022-Q3 / 2022-Q1 rows.
# ---------------------------------------------------------------------
$totals = $wb.Worksheets.Item(1)

# Duplicate the existing "index column" formatting (bold, thin border,
# centered) down into the new row positions before we overwrite values.
$totals.Range("A2").Copy($totals.Range("A3"))
$totals.Range("A2").Copy($totals.Range("A4"))

# Old row 3 (2022-Q1) -> new row 4, with index updated to 2
$totals.Range("A4").Value = 2
$totals.Range("B4").Value = "2022-Q1"
$totals.Range("C4").Value = 2
$totals.Range("D4").Value = 0

# Old row 2 (2022-Q3) -> new row 3, index stays 1
$totals.Range("A3").Value = 1
$totals.Range("B3").Value = "2022-Q3"
$totals.Range("C3").Value = 3
$totals.Range("D3").Value = 0.08

# New row 2: 2022-Q4, index 0
$totals.Range("A2").Value = 0
$totals.Range("B2").Value = "2022-Q4"
$totals.Range("C2").Value = 2
$totals.Range("D2").Value = 0.08

# ---------------------------------------------------------------------
# 2) Insert the new "2022-Q4" worksheet right after "总计".
# ---------------------------------------------------------------------
$q4 = $wb.Worksheets.Add($null, $totals)
$q4.Name = "2022-Q4"

# ---------------------------------------------------------------------
# 3) Populate "2022-Q4" with the fund-holdings table.
# ---------------------------------------------------------------------
$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

# NOTE: fund code / size / weight columns are stored as *text* in the
# source data (leading zeros, fixed trailing-zero decimals) - a leading
# apostrophe forces Excel to keep them as text instead of re-parsing them
# as numbers (which would drop the leading "0" / trailing zeros).
$q4.Range("A2").Value = 0
$q4.Range("B2").Value = "'016029"
$q4.Range("C2").Value = "湘财成长优选一年持有期混合A"
$q4.Range("D2").Value = "'1.70"
$q4.Range("E2").Value = "'89.45"
$q4.Range("F2").Value = "'4.59"
$q4.Range("G2").Value = "'0.0780"
$q4.Range("H2").Value = 2

$q4.Range("A3").Value = 1
$q4.Range("B3").Value = "'016030"
$q4.Range("C3").Value = "湘财成长优选一年持有期混合C"
$q4.Range("D3").Value = "'0.15"
$q4.Range("E3").Value = "'89.45"
$q4.Range("F3").Value = "'4.59"
$q4.Range("G3").Value = "'0.0069"
$q4.Range("H3").Value = 2

# Formatting to match the look of the sibling quarter sheets: bold,
# thin-bordered, centered header row and index column.
$headerRange = $q4.Range("B1:H1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

$indexRange = $q4.Range("A2:A3")
$indexRange.Font.Bold = $true
$indexRange.Borders.LineStyle = 1
$indexRange.HorizontalAlignment = -4108
$indexRange.VerticalAlignment = -4160
